$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 37038148
$ws.Range("I9").Value = 166667170
$ws.Range("J9").Value = 1285.2858
$ws.Range("K9").Value = 166667170
$ws.Range("L9").Value = 1285.2858
$ws.Range("M9").Value = -166667001
$ws.Range("N9").Value = -1623.2858

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 7386
$ws.Range("I80").Value = 324.2857
$ws.Range("K80").Value = 972.8571000000001
$ws.Range("M80").Value = 25.14289999999994

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 7386
$ws.Range("I83").Value = 324.2857
$ws.Range("K83").Value = 2918.5713
$ws.Range("M83").Value = 2073.4287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H114").Value = 81245
$ws.Range("J114").Value = 89994.42999999999
$ws.Range("L114").Value = 89994.42999999999
$ws.Range("N114").Value = -98672.42999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H120").Value = 47997.332
$ws.Range("J120").Value = 47997.332
$ws.Range("L120").Value = 47997.332
$ws.Range("N120").Value = -57673.332

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 76570.91
$ws.Range("J123").Value = 79228
$ws.Range("L123").Value = 79228
$ws.Range("N123").Value = -89028

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 364060.44
$ws.Range("I137").Value = 1173.0646
$ws.Range("J137").Value = 1614005.9
$ws.Range("K137").Value = 3519.1938
$ws.Range("L137").Value = 4842017.699999999
$ws.Range("M137").Value = -969.1938
$ws.Range("N137").Value = -4847117.699999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 47669324
$ws.Range("J138").Value = 76926280
$ws.Range("L138").Value = 230778840
$ws.Range("N138").Value = -230789120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 59567.8
$ws.Range("J7").Value = 60709.75
$ws.Range("L7").Value = 60709.75
$ws.Range("N7").Value = -60937.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 39932
$ws.Range("J104").Value = 39932
$ws.Range("L104").Value = 39932
$ws.Range("N104").Value = -46920

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 78228
$ws.Range("J107").Value = 78228
$ws.Range("L107").Value = 78228
$ws.Range("N107").Value = -85908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H108").Value = 88996.664
$ws.Range("J108").Value = 88996.664
$ws.Range("L108").Value = 88996.664
$ws.Range("N108").Value = -96676.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H117").Value = 53578.168
$ws.Range("J117").Value = 53578.168
$ws.Range("L117").Value = 53578.168
$ws.Range("N117").Value = -62756.168

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 92124.625
$ws.Range("J121").Value = 92124.625
$ws.Range("L121").Value = 92124.625
$ws.Range("N121").Value = -95618.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 50499.5
$ws.Range("J13").Value = 50499.5
$ws.Range("L13").Value = 50499.5
$ws.Range("N13").Value = -50835.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5996.3335
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 5996.3335
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 5996.3335
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -6490.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1929.0435
$ws.Range("I107").Value = 1806
$ws.Range("J107").Value = 2210.2856
$ws.Range("K107").Value = 1806
$ws.Range("L107").Value = 2210.2856
$ws.Range("M107").Value = 114
$ws.Range("N107").Value = -6050.2856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 76996.42999999999
$ws.Range("J109").Value = 76996.42999999999
$ws.Range("L109").Value = 76996.42999999999
$ws.Range("N109").Value = -79770.42999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H115").Value = 77911.836
$ws.Range("J115").Value = 77911.836
$ws.Range("L115").Value = 77911.836
$ws.Range("N115").Value = -81045.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 24516.572
$ws.Range("J132").Value = 24516.572
$ws.Range("L132").Value = 24516.572
$ws.Range("N132").Value = -34636.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 85494.875
$ws.Range("J138").Value = 85494.875
$ws.Range("L138").Value = 85494.875
$ws.Range("N138").Value = -95774.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 74988.57000000001
$ws.Range("J140").Value = 74988.57000000001
$ws.Range("L140").Value = 74988.57000000001
$ws.Range("N140").Value = -85348.57000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 957
$ws.Range("I14").Value = 913
$ws.Range("J14").Value = 1111
$ws.Range("K14").Value = 913
$ws.Range("L14").Value = 1111
$ws.Range("M14").Value = -743
$ws.Range("N14").Value = -1451

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 7897.625
$ws.Range("J15").Value = 12632.333
$ws.Range("L15").Value = 12632.333
$ws.Range("N15").Value = -12972.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 24975.334
$ws.Range("I18").Value = 24950
$ws.Range("K18").Value = 24950
$ws.Range("M18").Value = -24720

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 569.3333
$ws.Range("I25").Value = 708
$ws.Range("K25").Value = 708
$ws.Range("M25").Value = -534

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3252896.5
$ws.Range("I132").Value = 3794656.5
$ws.Range("J132").Value = 2169376.2
$ws.Range("K132").Value = 11383969.5
$ws.Range("L132").Value = 6508128.600000001
$ws.Range("M132").Value = -11381439.5
$ws.Range("N132").Value = -6513188.600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 6833.3335
$ws.Range("I76").Value = 2750
$ws.Range("K76").Value = 8250
$ws.Range("M76").Value = -7867

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H79").Value = 6833.3335
$ws.Range("I79").Value = 2750
$ws.Range("K79").Value = 8250
$ws.Range("M79").Value = -6924

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 4714.2856
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 4714.2856
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 14142.8568
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -19384.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 812.375
$ws.Range("I117").Value = 333.16666
$ws.Range("J117").Value = 2250
$ws.Range("K117").Value = 999.4999799999999
$ws.Range("L117").Value = 6750
$ws.Range("M117").Value = 2442.50002
$ws.Range("N117").Value = -13634

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 667956.6
$ws.Range("I121").Value = 1347.091
$ws.Range("K121").Value = 4041.273
$ws.Range("M121").Value = -2731.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1417.7059
$ws.Range("I131").Value = 958.4
$ws.Range("J131").Value = 2073.8572
$ws.Range("K131").Value = 2875.2
$ws.Range("L131").Value = 6221.571599999999
$ws.Range("M131").Value = 2164.8
$ws.Range("N131").Value = -16301.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 8824.25
$ws.Range("I99").Value = 8824.25
$ws.Range("K99").Value = 8824.25
$ws.Range("M99").Value = -6578.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 99990
$ws.Range("J108").Value = 99990
$ws.Range("L108").Value = 99990
$ws.Range("N108").Value = -107670

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 99757.60000000001
$ws.Range("J135").Value = 99757.60000000001
$ws.Range("L135").Value = 99757.60000000001
$ws.Range("N135").Value = -109897.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7627
$ws.Range("I61").Value = 7627
$ws.Range("K61").Value = 7627
$ws.Range("M61").Value = -7425

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4835
$ws.Range("I68").Value = 4801
$ws.Range("K68").Value = 4801
$ws.Range("M68").Value = -4052

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 4835
$ws.Range("I71").Value = 4801
$ws.Range("K71").Value = 24005
$ws.Range("M71").Value = -20261

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 7627
$ws.Range("I113").Value = 7627
$ws.Range("K113").Value = 7627
$ws.Range("M113").Value = -5457

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 45255.363
$ws.Range("J121").Value = 45255.363
$ws.Range("L121").Value = 45255.363
$ws.Range("N121").Value = -48749.363

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3863.5881
$ws.Range("I136").Value = 3570.077
$ws.Range("J136").Value = 4817.5
$ws.Range("K136").Value = 10710.231
$ws.Range("L136").Value = 14452.5
$ws.Range("M136").Value = -8160.231
$ws.Range("N136").Value = -19552.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 33378400
$ws.Range("I31").Value = 50030000
$ws.Range("J31").Value = 75200
$ws.Range("K31").Value = 50030000
$ws.Range("L31").Value = 75200
$ws.Range("M31").Value = -50029652
$ws.Range("N31").Value = -75896

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 61963.332
$ws.Range("J54").Value = 75445
$ws.Range("L54").Value = 75445
$ws.Range("N54").Value = -76485

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1195
$ws.Range("I107").Value = 1663.6666
$ws.Range("J107").Value = 811.5454999999999
$ws.Range("K107").Value = 4990.9998
$ws.Range("L107").Value = 2434.6365
$ws.Range("M107").Value = -3070.9998
$ws.Range("N107").Value = -6274.6365

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 47236.855
$ws.Range("J121").Value = 47236.855
$ws.Range("L121").Value = 47236.855
$ws.Range("N121").Value = -50730.855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3020.261
$ws.Range("I122").Value = 2986.3
$ws.Range("K122").Value = 8958.900000000001
$ws.Range("M122").Value = -6508.900000000001
